$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted at row 253. Inserting a whole row
# shifts every existing row (253..357) down by one - which reproduces the
# data shift seen in the diff (new row i = old row i-1 for i = 254..358) -
# and extends the sheet dimension from T357 to T358 automatically.
$ws.Rows.Item(253).Insert()

# The newly inserted row 253 is completely blank. Re-use the row that now
# sits right below it (row 254, which holds the data that used to be in
# row 253) as a template, then overwrite just the cells that actually
# differ for the new observation.
$ws.Rows.Item(254).Copy()
$ws.Rows.Item(253).PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("D253").Value = 44609
$ws.Range("M253").Value = 200
$ws.Range("N253").Value = 17000
$ws.Range("O253").Value = 17000
$ws.Range("P253").Value = 17000
$ws.Range("S253").Value = 1133
